# Trading update: 2026-02-17 23:54:45
# Append the new trade row (row 21) to the "All Trades" and "MarketMaking"
# sheets - both sheets keep an identical trade log.

$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "MarketMaking")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = 21

    $ws.Cells.Item($row, 1).Value = 20
    # Leading apostrophe forces text storage so the date-shaped string isn't
    # reinterpreted as a date serial number (matches the source data, which
    # stores dates/times as plain text).
    $ws.Cells.Item($row, 2).Value = "'2026-02-17"
    $ws.Cells.Item($row, 3).Value = "23:54:39"
    $ws.Cells.Item($row, 4).Value = "MarketMaking"
    $ws.Cells.Item($row, 5).Value = "DOWN"
    $ws.Cells.Item($row, 6).Value = 0.01
    $ws.Cells.Item($row, 7).Value = ""
    $ws.Cells.Item($row, 8).Value = "OPEN"
    $ws.Cells.Item($row, 9).Value = 0
    $ws.Cells.Item($row, 10).Value = 0
    $ws.Cells.Item($row, 11).Value = 100.3236569789373
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($row, 16).Value = ""
    $ws.Cells.Item($row, 17).Value = 0
}
